$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.378.19"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.17"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.00"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6297"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07587"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2929"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.46"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.840.67"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.002"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001082"
$ws.Range("E14").Value = "  +7.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6781"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.63"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.089.75"
$ws.Range("E17").Value = "  -7.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.159"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.406.49"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.70"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.419"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.95"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.388"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.312"
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.461"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05603"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.033"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7098"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.582"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.233.19"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.454"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9072"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.64"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.02"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.206"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.977"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.680"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1121"
$ws.Range("E51").Value = "  -0.65%  "
